$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = "43% CR/LFM+CDM/H:1`n2% CR/LFM+CDM/H:1`n35% CR+PC/LFM+CDL/H:1`n5% S+SL/LFM+CDM/H:1`n11% S/LFM+CDM/H:1`n3% S/LFBR+CDH/H:1`n1% S/LFBR+CDM/H:1"

$ws.Range("B5").Select()
